# Updated cryptos list on Wed Oct 18 10:50:25 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.456.29'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.26%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.584.29'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.20%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.32%  '

$ws.Range("E6").Value = '  +0.13%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.36'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '24.06'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.90%  '

$ws.Range("E10").Value = '  -1.96%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0592'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.45%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0894'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.92%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.809.78'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.23%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.583.97'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.33%  '

$ws.Range("E15").Value = '  -0.88%  '

$ws.Range("E16").Value = '  -1.50%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '28.480.55'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.11%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.20'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.43%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '230.89'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.31%  '

$ws.Range("E20").Value = '  -0.44%  '

$ws.Range("E21").Value = '  -2.25%  '

$ws.Range("E22").Value = '  +0.10%  '

$ws.Range("E23").Value = '  -3.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.85%  '

$ws.Range("E25").Value = '  +3.96%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.12'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.26%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.04'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.36%  '

$ws.Range("E28").Value = '  -1.56%  '

$ws.Range("E29").Value = '  -1.93%  '

$ws.Range("E30").Value = '  +0.06%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0482'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.46%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.22'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.00%  '

$ws.Range("E34").Value = '  -2.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.399.96'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.04%  '

$ws.Range("E36").Value = '  +6.19%  '

$ws.Range("E37").Value = '  -4.22%  '

$ws.Range("E38").Value = '  +0.48%  '

$ws.Range("E39").Value = '  +1.86%  '

$ws.Range("E40").Value = '  -0.82%  '

$ws.Range("E41").Value = '  -3.62%  '

$ws.Range("E42").Value = '  +0.08%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.792'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.58%  '

$ws.Range("E44").Value = '  +1.29%  '

$ws.Range("E45").Value = '  -0.62%  '

$ws.Range("E46").Value = '  -3.25%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.961'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.30%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '63.29'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.46%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.721.82'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.14%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '86.75'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.61%  '

$ws.Range("E51").Value = '  -0.89%  '
